$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-design the header row (no index column change, just rename headers)
$ws.Range("A1").Value = "Nome da Localidade"
$ws.Range("B1").Value = "LONG"
$ws.Range("C1").Value = "LAT"
$ws.Range("D1").Value = "Percentual"
$ws.Range("E1").Value = "Total pesquia"
$ws.Range("F1").Value = "Votos"

# Widen column A to fit new header text (best-fit autosize like Excel does)
$ws.Columns.Item(1).EntireColumn.AutoFit()

# Update selection
$ws.Range("K7").Select()
